{"js": "// Apply the LOQ4054.docx edits described in the commit diff.\n// Each change is a full-text replacement of a single <w:t> run's\n// content, located via Body.search() (exact, case-sensitive match)\n// and applied with Range.insertText(..., \"Replace\").\nconst body = context.document.body;\n\nconst replacements = [\n  { find: \"Ativa\u00e7\u00e3o: 01/01/2018\", replace: \"Ativa\u00e7\u00e3o: 01/01/2024\" },\n  { find: \"Curso (semestre ideal): EB (7), EQD (6), EQN (7)\", replace: \"Curso (semestre ideal): EB (7), EQN (7)\" },\n  { find: \"1- Introdu\u00e7\u00e3o: Transfer\u00eancia de massa: Defini\u00e7\u00e3o. Classifica\u00e7\u00e3o das opera\u00e7\u00f5es que envolvem transfer\u00eancia de massa. Contribui\u00e7\u00f5es \u00e0 transfer\u00eancia de massa. Tipos de difus\u00e3o.2- Coeficiente de difus\u00e3o: Considera\u00e7\u00f5es. Difus\u00e3o em gases: An\u00e1lise da primeira lei de Fick; O coeficiente de difus\u00e3o para gases. Estimativa do coeficiente de difus\u00e3o a partir de um coeficiente de difus\u00e3o conhecido em alta temperatura e press\u00e3o. Coeficiente de difus\u00e3o de um soluto em uma mistura gasosa estagnada de multicomponentes. Difus\u00e3o em l\u00edquidos. Difus\u00e3o em s\u00f3lidos.3- Concentra\u00e7\u00f5es, velocidade e fluxos: Concentra\u00e7\u00e3o. Velocidade. Fluxo. A equa\u00e7\u00e3o de Stefan \u2013 Maxwel.4 - Equa\u00e7\u00f5es da continuidade em transfer\u00eancia de massa: Equa\u00e7\u00f5es da continuidade molar de um soluto. Regime transiente sem/com velocidade do meio nula. Meio sem e com rea\u00e7\u00e3o qu\u00edmica.5- Difus\u00e3o em regime permanente sem rea\u00e7\u00e3o qu\u00edmica: Difus\u00e3o Unidimensional em regime permanente. Difus\u00e3o atrav\u00e9s de filme gasoso inerte e estagnado. Difus\u00e3o pseudo-estacion\u00e1ria num filme gasoso estagnado. Contradifus\u00e3o equimolar. Taxa molar em esferas isoladas. Difus\u00e3o em membranas.6- Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica: Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea na superf\u00edcie de uma part\u00edcula catal\u00edtica n\u00e3o porosa. Difus\u00e3o com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea na superf\u00edcie de uma part\u00edcula n\u00e3o catal\u00edtica e n\u00e3o porosa. Difus\u00e3o intrapart\u00edcula com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea. Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica homog\u00eanea.7- Transfer\u00eancia de massa entre fases: Teoria das duas resist\u00eancias. Coeficiente individual e global de transfer\u00eancia de massa. Coeficientes globais de transfer\u00eancia de massa. Coeficientes volum\u00e9tricos de transfer\u00eancia de massa para torres de recheios. Balan\u00e7o macrosc\u00f3pio de mat\u00e9ria. Opera\u00e7\u00f5es cont\u00ednuas.\", replace: \"1- Introdu\u00e7\u00e3o: Transfer\u00eancia de massa: Defini\u00e7\u00e3o. Classifica\u00e7\u00e3o das opera\u00e7\u00f5es que envolvem transfer\u00eancia de massa. Contribui\u00e7\u00f5es \u00e0 transfer\u00eancia de massa. Tipos de difus\u00e3o. 2- Coeficiente e mecanismos de difus\u00e3o: Considera\u00e7\u00f5es a respeito. Difus\u00e3o em gases: An\u00e1lise da primeira lei de Fick. O coeficiente de difus\u00e3o para gases. Estimativa do coeficiente de difus\u00e3o a partir de um coeficiente de difus\u00e3o conhecido em outra temperatura e press\u00e3o. Coeficiente de difus\u00e3o de um soluto em uma mistura gasosa estagnada de multicomponentes. Difus\u00e3o em l\u00edquidos. Difus\u00e3o em s\u00f3lidos. 3- Concentra\u00e7\u00f5es, velocidades e fluxos: Concentra\u00e7\u00e3o. Velocidade. Fluxo. A equa\u00e7\u00e3o de Stefan \u2013 Maxwel. Coeficiente convectivo de transfer\u00eancia de massa 4 - Equa\u00e7\u00f5es da continuidade em transfer\u00eancia de massa: Considera\u00e7\u00f5es a respeito. Equa\u00e7\u00f5es da continuidade m\u00e1ssica e molar de um soluto. Equa\u00e7\u00f5es da continuidade do soluto A em termos da lei ordin\u00e1ria da difus\u00e3o. Condi\u00e7\u00f5es de contorno. 5- Difus\u00e3o em regime permanente sem rea\u00e7\u00e3o qu\u00edmica: Difus\u00e3o Unidimensional em regime permanente. Difus\u00e3o atrav\u00e9s de filme gasoso inerte e estagnado. Difus\u00e3o pseudo-estacion\u00e1ria num filme gasoso estagnado. Contradifus\u00e3o equimolar. Taxa molar em esferas isoladas. Difus\u00e3o em membranas. 6- Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica: Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea na superf\u00edcie de uma part\u00edcula catal\u00edtica n\u00e3o porosa. Difus\u00e3o com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea na superf\u00edcie de uma part\u00edcula n\u00e3o catal\u00edtica e n\u00e3o porosa. Difus\u00e3o intraparticular com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea. Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica homog\u00eanea. 7- Transfer\u00eancia de massa entre fases: Considera\u00e7\u00f5es a respeito. T\u00e9cnicas de separa\u00e7\u00e3o. Transfer\u00eancia de massa entre fases. Teoria das duas resist\u00eancias. Coeficientes globais de transfer\u00eancia de massa. Coeficientes volum\u00e9tricos de transfer\u00eancia de massa para torres de recheios. Balan\u00e7o macrosc\u00f3pio de mat\u00e9ria em equipamentos de separa\u00e7\u00e3o. Opera\u00e7\u00f5es cont\u00ednuas (contracorrente e paralelo). C\u00e1lculo da altura efetiva e do di\u00e2metro de uma coluna para opera\u00e7\u00e3o cont\u00ednua em um sistema dilu\u00eddo.\" },\n  { find: \"1 - Introduction: Mass transfer: Definition. Classification of operations involving mass transfer. Contributions to mass transfer. Types of diffusion. 2 - Diffusion coefficient: Considerations. Diffusion in gases: Analysis of Fick's first law, the diffusion coefficient for gases. Estimation of the diffusion coefficient from a known diffusion coefficient at high temperature and pressure. Diffusion coefficient of a solute in a multicomponent mixture of stagnant gases. Diffusion in liquids. Diffusion in solids. 3 - Concentrations, and flow rate: Concentration. Speed and Flow. The equation of Stefan \u2013 Maxwell. 4 - Equation of continuity for mass transfer: Equations of continuity of a solute molar. Transient without/with speed zero means. Medium with and without chemical reaction. 5 - Diffusion in permanent regime without chemical reaction: One-dimensional diffusion in permanent regime. Diffusion through inert and stagnant gaseous film. Pseudo-stationary diffusion in a stagnant gaseous film. Equimolar contradifusion. Molar rate in isolated beads. Diffusion in membranes. 6 - Diffusion in permanent with chemical reaction: Diffusion in steady state with heterogeneous chemical reaction on the surface of a nonporous catalytic particle. Diffusion with heterogeneous chemical reaction on the surface of a non-catalytic and non-porous particle. Intraparticle diffusion with heterogeneous chemical reaction. Continuous diffusion with homogeneous chemical reaction. 7 - Mass transfer between phases: Theory of the two resistors. Individual and global coefficient of mass transfer. Global mass transfer coefficients. Volumetric coefficients of transfer of mass to towers of fillings. Balance macroscope of matter. Continuous operations.\", replace: \"1- Introduction: Mass transfer: Definition. Classification of operations involving mass transfer. Contributions to mass transfer. Diffusion types. 2- Coefficient and diffusion mechanisms: Considerations in this regard. Diffusion in gases: Analysis of Fick's first law. The diffusion coefficient for gases. Estimation of the diffusion coefficient from a known diffusion coefficient at another temperature and pressure. Diffusion coefficient of a solute in a stagnant multicomponent gaseous mixture. Diffusion in liquids. Diffusion in solids. 3- Concentrations, velocities, and flows: Concentration. Speed. Flow. The Stefan\u2013Maxwell equation. Convective coefficient of mass transfer 4 - Equations of continuity in mass transfer: Considerations in this regard. Mass and molar continuity equations for a solute. Continuity equations for solute A in terms of the ordinary law of diffusion. Boundary conditions. 5- Steady state diffusion without chemical reaction: Unidimensional steady state diffusion. Diffusion through inert and stagnant gaseous film. Pseudo-stationary diffusion in a stagnant gas film. Equimolar counterdiffusion. Molar rate in isolated spheres. Diffusion in membranes. 6- Steady-state diffusion with chemical reaction: Steady-state diffusion with heterogeneous chemical reaction on the surface of a non-porous catalytic particle. Diffusion with heterogeneous chemical reaction on the surface of a non-catalytic, non-porous particle. Intraparticulate diffusion with heterogeneous chemical reaction. Diffusion in steady state with homogeneous chemical reaction. 7- Mass transfer between phases: Considerations in this regard. Separation techniques. Mass transfer between phases. Theory of two resistances. Global mass transfer coefficients. Volumetric mass transfer coefficients for packing towers. Macroscopic balance of matter in separation equipment. Continuous operations (countercurrent and parallel). Calculation of the effective height and diameter of a column for continuous operation in a dilute system.\" },\n  { find: \"A avalia\u00e7\u00e3o ser\u00e1 feita por meio de provas escritas.\", replace: \"A avalia\u00e7\u00e3o ser\u00e1 feita por meio de provas escritas e trabalhos em grupos.\" },\n  { find: \"A Nota Final (NF) ser\u00e1 calculada da seguinte maneira: NF = (P1 + 2*P2)/3.\", replace: \"A Nota Final (NF) ser\u00e1 calculada da seguinte maneira: NF = (P1 + 2*P2)/3P2 = Nota da Prova (80%) e Nota do Trabalho (20%).\" },\n  { find: \"1) CREMASCO, M. A. Fundamentos de Transfer\u00eancia de Massa. Campinas: Editora Unicamp; 2008-2009.2) INCROPERA, F. P.; WITT, D. P. Fundamentos de Transfer\u00eancia de Calor e Massa. Rio de Janeiro: LTC, 2008.3) Bird, R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fen\u00f4menos de Transporte. 2\u00aa ed. Rio de Janeiro: LTC, 2004.4) BENNETT, C. O.; MYERS, J. E. Fen\u00f4meno de Transporte: Quantidade de Movimento, Calor e Massa. S\u00e3o Paulo: McGrawc- Hill, 1978.5) COULSON, J. M.; RICHARDSON, J. F.; BACKHURST, J. R.; HARKER, J. H. Fluid Flow, Heat Transfer and Mass Transfer. In: COULSON & Richardson Series - Chemical Engineering. 5th ed. Pergamon Press, Oxford, 1996. v.16) FOUST, A. S.; Wenzel, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSON, L. B. Princ\u00edpios das Opera\u00e7\u00f5es Unit\u00e1rias. 2\u00aa ed. Rio de Janeiro: Guanabara Dois, 1982.7) PERRY's Chemical Engineers Handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry. New York: McGraw-Hill, 2008.8) WELTY, J. R.; PIGFORD, R. L.; WILKE, C. R. Fundamentals of Momentum, Heat, and Mass Transfer. 5th ed. USA: John Wiley & Sons, Inc, 2008.9) POLING, B. E.; PRAUSNITZ, J. M.; O'CONNELL, J. The Properties of Gases and Liquids. 5th ed. New York: McGraw-Hill, 2000.10) CALDAS, J. N.; DE LACERDA, A. I.; VELOSO, E.; PASCHOAL, L. C. M. Internos de Torres: Pratos & Recheios. 2\u00aa ed. Rio de Janeiro: Editora Interci\u00eancia, 2007.\", replace: \"1) CREMASCO, M. A. Fundamentos de Transfer\u00eancia de Massa, 3\u00aa ed. S\u00e3o Paulo: Editora Blucher; 2021. 2) INCROPERA, F. P.; WITT, D. P. Fundamentos de Transfer\u00eancia de Calor e Massa. 8\u00aa ed. Rio de Janeiro: LTC, 2019. 3) Bird, R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fen\u00f4menos de Transporte. 2\u00aa ed. Rio de Janeiro: LTC, 2004. 4) COULSON, J. M.; RICHARDSON, J. F.; BACKHURST, J. R.; HARKER, J. H. Fluid Flow, Heat Transfer and Mass Transfer. In: COULSON & Richardson Series - Chemical Engineering. 6th ed. Pergamon Press, Oxford, 1999. v.1 5) PERRY's Chemical Engineers Handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry. 9\u00aa ed. New York: McGraw-Hill, 2019. 6) WELTY, J. R.; PIGFORD, R. L.; WILKE, C. R. Fundamentals of Momentum, Heat, and Mass Transfer. 6th ed. USA: John Wiley & Sons, Inc, 2014. 7) POLING, B. E.; PRAUSNITZ, J. M.; O'CONNELL, J. The Properties of Gases and Liquids. 5th ed. New York: McGraw-Hill, 2004. 8) CALDAS, J. N.; DE LACERDA, A. I.; VELOSO, E.; PASCHOAL, L. C. M. Internos de Torres: Pratos & Recheios. 2\u00aa ed. Rio de Janeiro: Editora Interci\u00eancia, 2007.\" },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + find.substring(0, 60));\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Apply the LOQ4054.docx edits described in the commit diff.\n# Each change is a full-text replacement: Find.Execute locates the\n# exact (case-sensitive, non-wildcard) run text across the whole\n# document body, then the matched Range's .Text is overwritten\n# directly (Range.Text=, not Find's ReplaceWith, so Word's\n# Find/Replace \"smart quotes\" AutoFormat never touches straight\n# apostrophes inside the new text).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = 'Ativa\u00e7\u00e3o: 01/01/2018'; Replace = 'Ativa\u00e7\u00e3o: 01/01/2024' },\n    @{ Find = 'Curso (semestre ideal): EB (7), EQD (6), EQN (7)'; Replace = 'Curso (semestre ideal): EB (7), EQN (7)' },\n    @{ Find = '1- Introdu\u00e7\u00e3o: Transfer\u00eancia de massa: Defini\u00e7\u00e3o. Classifica\u00e7\u00e3o das opera\u00e7\u00f5es que envolvem transfer\u00eancia de massa. Contribui\u00e7\u00f5es \u00e0 transfer\u00eancia de massa. Tipos de difus\u00e3o.2- Coeficiente de difus\u00e3o: Considera\u00e7\u00f5es. Difus\u00e3o em gases: An\u00e1lise da primeira lei de Fick; O coeficiente de difus\u00e3o para gases. Estimativa do coeficiente de difus\u00e3o a partir de um coeficiente de difus\u00e3o conhecido em alta temperatura e press\u00e3o. Coeficiente de difus\u00e3o de um soluto em uma mistura gasosa estagnada de multicomponentes. Difus\u00e3o em l\u00edquidos. Difus\u00e3o em s\u00f3lidos.3- Concentra\u00e7\u00f5es, velocidade e fluxos: Concentra\u00e7\u00e3o. Velocidade. Fluxo. A equa\u00e7\u00e3o de Stefan \u2013 Maxwel.4 - Equa\u00e7\u00f5es da continuidade em transfer\u00eancia de massa: Equa\u00e7\u00f5es da continuidade molar de um soluto. Regime transiente sem/com velocidade do meio nula. Meio sem e com rea\u00e7\u00e3o qu\u00edmica.5- Difus\u00e3o em regime permanente sem rea\u00e7\u00e3o qu\u00edmica: Difus\u00e3o Unidimensional em regime permanente. Difus\u00e3o atrav\u00e9s de filme gasoso inerte e estagnado. Difus\u00e3o pseudo-estacion\u00e1ria num filme gasoso estagnado. Contradifus\u00e3o equimolar. Taxa molar em esferas isoladas. Difus\u00e3o em membranas.6- Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica: Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea na superf\u00edcie de uma part\u00edcula catal\u00edtica n\u00e3o porosa. Difus\u00e3o com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea na superf\u00edcie de uma part\u00edcula n\u00e3o catal\u00edtica e n\u00e3o porosa. Difus\u00e3o intrapart\u00edcula com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea. Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica homog\u00eanea.7- Transfer\u00eancia de massa entre fases: Teoria das duas resist\u00eancias. Coeficiente individual e global de transfer\u00eancia de massa. Coeficientes globais de transfer\u00eancia de massa. Coeficientes volum\u00e9tricos de transfer\u00eancia de massa para torres de recheios. Balan\u00e7o macrosc\u00f3pio de mat\u00e9ria. Opera\u00e7\u00f5es cont\u00ednuas.'; Replace = '1- Introdu\u00e7\u00e3o: Transfer\u00eancia de massa: Defini\u00e7\u00e3o. Classifica\u00e7\u00e3o das opera\u00e7\u00f5es que envolvem transfer\u00eancia de massa. Contribui\u00e7\u00f5es \u00e0 transfer\u00eancia de massa. Tipos de difus\u00e3o. 2- Coeficiente e mecanismos de difus\u00e3o: Considera\u00e7\u00f5es a respeito. Difus\u00e3o em gases: An\u00e1lise da primeira lei de Fick. O coeficiente de difus\u00e3o para gases. Estimativa do coeficiente de difus\u00e3o a partir de um coeficiente de difus\u00e3o conhecido em outra temperatura e press\u00e3o. Coeficiente de difus\u00e3o de um soluto em uma mistura gasosa estagnada de multicomponentes. Difus\u00e3o em l\u00edquidos. Difus\u00e3o em s\u00f3lidos. 3- Concentra\u00e7\u00f5es, velocidades e fluxos: Concentra\u00e7\u00e3o. Velocidade. Fluxo. A equa\u00e7\u00e3o de Stefan \u2013 Maxwel. Coeficiente convectivo de transfer\u00eancia de massa 4 - Equa\u00e7\u00f5es da continuidade em transfer\u00eancia de massa: Considera\u00e7\u00f5es a respeito. Equa\u00e7\u00f5es da continuidade m\u00e1ssica e molar de um soluto. Equa\u00e7\u00f5es da continuidade do soluto A em termos da lei ordin\u00e1ria da difus\u00e3o. Condi\u00e7\u00f5es de contorno. 5- Difus\u00e3o em regime permanente sem rea\u00e7\u00e3o qu\u00edmica: Difus\u00e3o Unidimensional em regime permanente. Difus\u00e3o atrav\u00e9s de filme gasoso inerte e estagnado. Difus\u00e3o pseudo-estacion\u00e1ria num filme gasoso estagnado. Contradifus\u00e3o equimolar. Taxa molar em esferas isoladas. Difus\u00e3o em membranas. 6- Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica: Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea na superf\u00edcie de uma part\u00edcula catal\u00edtica n\u00e3o porosa. Difus\u00e3o com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea na superf\u00edcie de uma part\u00edcula n\u00e3o catal\u00edtica e n\u00e3o porosa. Difus\u00e3o intraparticular com rea\u00e7\u00e3o qu\u00edmica heterog\u00eanea. Difus\u00e3o em regime permanente com rea\u00e7\u00e3o qu\u00edmica homog\u00eanea. 7- Transfer\u00eancia de massa entre fases: Considera\u00e7\u00f5es a respeito. T\u00e9cnicas de separa\u00e7\u00e3o. Transfer\u00eancia de massa entre fases. Teoria das duas resist\u00eancias. Coeficientes globais de transfer\u00eancia de massa. Coeficientes volum\u00e9tricos de transfer\u00eancia de massa para torres de recheios. Balan\u00e7o macrosc\u00f3pio de mat\u00e9ria em equipamentos de separa\u00e7\u00e3o. Opera\u00e7\u00f5es cont\u00ednuas (contracorrente e paralelo). C\u00e1lculo da altura efetiva e do di\u00e2metro de uma coluna para opera\u00e7\u00e3o cont\u00ednua em um sistema dilu\u00eddo.' },\n    @{ Find = '1 - Introduction: Mass transfer: Definition. Classification of operations involving mass transfer. Contributions to mass transfer. Types of diffusion. 2 - Diffusion coefficient: Considerations. Diffusion in gases: Analysis of Fick''s first law, the diffusion coefficient for gases. Estimation of the diffusion coefficient from a known diffusion coefficient at high temperature and pressure. Diffusion coefficient of a solute in a multicomponent mixture of stagnant gases. Diffusion in liquids. Diffusion in solids. 3 - Concentrations, and flow rate: Concentration. Speed and Flow. The equation of Stefan \u2013 Maxwell. 4 - Equation of continuity for mass transfer: Equations of continuity of a solute molar. Transient without/with speed zero means. Medium with and without chemical reaction. 5 - Diffusion in permanent regime without chemical reaction: One-dimensional diffusion in permanent regime. Diffusion through inert and stagnant gaseous film. Pseudo-stationary diffusion in a stagnant gaseous film. Equimolar contradifusion. Molar rate in isolated beads. Diffusion in membranes. 6 - Diffusion in permanent with chemical reaction: Diffusion in steady state with heterogeneous chemical reaction on the surface of a nonporous catalytic particle. Diffusion with heterogeneous chemical reaction on the surface of a non-catalytic and non-porous particle. Intraparticle diffusion with heterogeneous chemical reaction. Continuous diffusion with homogeneous chemical reaction. 7 - Mass transfer between phases: Theory of the two resistors. Individual and global coefficient of mass transfer. Global mass transfer coefficients. Volumetric coefficients of transfer of mass to towers of fillings. Balance macroscope of matter. Continuous operations.'; Replace = '1- Introduction: Mass transfer: Definition. Classification of operations involving mass transfer. Contributions to mass transfer. Diffusion types. 2- Coefficient and diffusion mechanisms: Considerations in this regard. Diffusion in gases: Analysis of Fick''s first law. The diffusion coefficient for gases. Estimation of the diffusion coefficient from a known diffusion coefficient at another temperature and pressure. Diffusion coefficient of a solute in a stagnant multicomponent gaseous mixture. Diffusion in liquids. Diffusion in solids. 3- Concentrations, velocities, and flows: Concentration. Speed. Flow. The Stefan\u2013Maxwell equation. Convective coefficient of mass transfer 4 - Equations of continuity in mass transfer: Considerations in this regard. Mass and molar continuity equations for a solute. Continuity equations for solute A in terms of the ordinary law of diffusion. Boundary conditions. 5- Steady state diffusion without chemical reaction: Unidimensional steady state diffusion. Diffusion through inert and stagnant gaseous film. Pseudo-stationary diffusion in a stagnant gas film. Equimolar counterdiffusion. Molar rate in isolated spheres. Diffusion in membranes. 6- Steady-state diffusion with chemical reaction: Steady-state diffusion with heterogeneous chemical reaction on the surface of a non-porous catalytic particle. Diffusion with heterogeneous chemical reaction on the surface of a non-catalytic, non-porous particle. Intraparticulate diffusion with heterogeneous chemical reaction. Diffusion in steady state with homogeneous chemical reaction. 7- Mass transfer between phases: Considerations in this regard. Separation techniques. Mass transfer between phases. Theory of two resistances. Global mass transfer coefficients. Volumetric mass transfer coefficients for packing towers. Macroscopic balance of matter in separation equipment. Continuous operations (countercurrent and parallel). Calculation of the effective height and diameter of a column for continuous operation in a dilute system.' },\n    @{ Find = 'A avalia\u00e7\u00e3o ser\u00e1 feita por meio de provas escritas.'; Replace = 'A avalia\u00e7\u00e3o ser\u00e1 feita por meio de provas escritas e trabalhos em grupos.' },\n    @{ Find = 'A Nota Final (NF) ser\u00e1 calculada da seguinte maneira: NF = (P1 + 2*P2)/3.'; Replace = 'A Nota Final (NF) ser\u00e1 calculada da seguinte maneira: NF = (P1 + 2*P2)/3P2 = Nota da Prova (80%) e Nota do Trabalho (20%).' },\n    @{ Find = '1) CREMASCO, M. A. Fundamentos de Transfer\u00eancia de Massa. Campinas: Editora Unicamp; 2008-2009.2) INCROPERA, F. P.; WITT, D. P. Fundamentos de Transfer\u00eancia de Calor e Massa. Rio de Janeiro: LTC, 2008.3) Bird, R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fen\u00f4menos de Transporte. 2\u00aa ed. Rio de Janeiro: LTC, 2004.4) BENNETT, C. O.; MYERS, J. E. Fen\u00f4meno de Transporte: Quantidade de Movimento, Calor e Massa. S\u00e3o Paulo: McGrawc- Hill, 1978.5) COULSON, J. M.; RICHARDSON, J. F.; BACKHURST, J. R.; HARKER, J. H. Fluid Flow, Heat Transfer and Mass Transfer. In: COULSON & Richardson Series - Chemical Engineering. 5th ed. Pergamon Press, Oxford, 1996. v.16) FOUST, A. S.; Wenzel, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSON, L. B. Princ\u00edpios das Opera\u00e7\u00f5es Unit\u00e1rias. 2\u00aa ed. Rio de Janeiro: Guanabara Dois, 1982.7) PERRY''s Chemical Engineers Handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry. New York: McGraw-Hill, 2008.8) WELTY, J. R.; PIGFORD, R. L.; WILKE, C. R. Fundamentals of Momentum, Heat, and Mass Transfer. 5th ed. USA: John Wiley & Sons, Inc, 2008.9) POLING, B. E.; PRAUSNITZ, J. M.; O''CONNELL, J. The Properties of Gases and Liquids. 5th ed. New York: McGraw-Hill, 2000.10) CALDAS, J. N.; DE LACERDA, A. I.; VELOSO, E.; PASCHOAL, L. C. M. Internos de Torres: Pratos & Recheios. 2\u00aa ed. Rio de Janeiro: Editora Interci\u00eancia, 2007.'; Replace = '1) CREMASCO, M. A. Fundamentos de Transfer\u00eancia de Massa, 3\u00aa ed. S\u00e3o Paulo: Editora Blucher; 2021. 2) INCROPERA, F. P.; WITT, D. P. Fundamentos de Transfer\u00eancia de Calor e Massa. 8\u00aa ed. Rio de Janeiro: LTC, 2019. 3) Bird, R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fen\u00f4menos de Transporte. 2\u00aa ed. Rio de Janeiro: LTC, 2004. 4) COULSON, J. M.; RICHARDSON, J. F.; BACKHURST, J. R.; HARKER, J. H. Fluid Flow, Heat Transfer and Mass Transfer. In: COULSON & Richardson Series - Chemical Engineering. 6th ed. Pergamon Press, Oxford, 1999. v.1 5) PERRY''s Chemical Engineers Handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry. 9\u00aa ed. New York: McGraw-Hill, 2019. 6) WELTY, J. R.; PIGFORD, R. L.; WILKE, C. R. Fundamentals of Momentum, Heat, and Mass Transfer. 6th ed. USA: John Wiley & Sons, Inc, 2014. 7) POLING, B. E.; PRAUSNITZ, J. M.; O''CONNELL, J. The Properties of Gases and Liquids. 5th ed. New York: McGraw-Hill, 2004. 8) CALDAS, J. N.; DE LACERDA, A. I.; VELOSO, E.; PASCHOAL, L. C. M. Internos de Torres: Pratos & Recheios. 2\u00aa ed. Rio de Janeiro: Editora Interci\u00eancia, 2007.' },\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $found = $rng.Find.Execute($r.Find)\n    if (-not $found) {\n        throw \"No match found for: \" + $r.Find.Substring(0, [Math]::Min(60, $r.Find.Length))\n    }\n    $rng.Text = $r.Replace\n}\n"}
